$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "EMPRESA"
$ws.Range("B1").Value = "RE"
$ws.Range("C1").Value = "NOME"
$ws.Range("D1").Value = "STATUS"

# Row 2 - NOME/STATUS entered before EMPRESA/RE to match shared string order
$ws.Range("C2").Value = "LUCAS"
$ws.Range("D2").Value = "ATIVO"
$ws.Range("A2").Value = "CARRARINHA"
$ws.Range("B2").Value = 94013

# Row 3
$ws.Range("A3").Value = "Safira"
$ws.Range("B3").Value = 12345
$ws.Range("C3").Value = "GUSTAVO"
$ws.Range("D3").Value = "ATIVO"

# Row 4
$ws.Range("A4").Value = "esf"
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = "DADÁ"
$ws.Range("D4").Value = "ATIVO"

$ws.Range("A5").Select()
